$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '29.110.38', '  +0.70%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.836.35', '  +0.74%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.001', '  +0.65%  ')
    ,@(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '242.42', '  -0.24%  ')
    ,@(6, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.6201', '  -1.47%  ')
    ,@(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.003', '  +0.73%  ')
    ,@(8, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07461', '  +0.17%  ')
    ,@(9, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2933', '  +0.11%  ')
    ,@(10, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '23.12', '  +0.60%  ')
    ,@(11, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07684', '  +0.13%  ')
    ,@(12, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.833.53', '  +0.64%  ')
    ,@(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.017', '  +0.87%  ')
    ,@(14, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6750', '  +1.46%  ')
    ,@(15, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '83.15', '  +0.44%  ')
    ,@(16, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000009142', '  -5.56%  ')
    ,@(17, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.923', '  -1.45%  ')
    ,@(18, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '29.096.50', '  +0.62%  ')
    ,@(19, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.084.60', '  +0.57%  ')
    ,@(20, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '241.47', '  +7.40%  ')
    ,@(21, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '12.73', '  +1.70%  ')
    ,@(22, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.003', '  +0.83%  ')
    ,@(23, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.213', '  +1.54%  ')
    ,@(24, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.002', '  +0.44%  ')
    ,@(25, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '159.15', '  -0.52%  ')
    ,@(26, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1417', '  +0.58%  ')
    ,@(27, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.519', '  +0.48%  ')
    ,@(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '17.93', '  +0.48%  ')
    ,@(29, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.501', '  +0.44%  ')
    ,@(30, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05614', '  +3.18%  ')
    ,@(31, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.125', '  +2.06%  ')
    ,@(32, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.140', '  +0.79%  ')
    ,@(33, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.206', '  +0.81%  ')
    ,@(34, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.846', '  -0.31%  ')
    ,@(35, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7429', '  -0.01%  ')
    ,@(36, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.144', '  +1.08%  ')
    ,@(37, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.661', '  +2.26%  ')
    ,@(38, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.780', '  +1.68%  ')
    ,@(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01788', '  +0.87%  ')
    ,@(40, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.212.86', '  -2.00%  ')
    ,@(41, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.401', '  -4.03%  ')
    ,@(42, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.8968', '  -0.19%  ')
    ,@(43, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '1.002', '  +0.69%  ')
    ,@(44, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '101.56', '  +0.48%  ')
    ,@(45, 'RocketPoolETH', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', '1.983.24', '  +0.60%  ')
    ,@(46, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '65.59', '  +1.10%  ')
    ,@(47, 'BabyDogeCoin', 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge', '0.00000000122', '  -1.13%  ')
    ,@(48, 'Mantle', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt', '0.5098', '  +0.63%  ')
    ,@(49, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4071', '  +0.86%  ')
    ,@(50, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.156', '  +2.17%  ')
    ,@(51, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05824', '  +0.75%  ')
)

foreach ($item in $data) {
    $row = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    # Prefix Price/Volume values with an apostrophe so Excel stores them as
    # literal text (matching the source inlineStr cells) instead of coercing
    # numeric-looking strings (e.g. '1.001') into floating point numbers.
    $ws.Cells.Item($row, 4).Value = "'" + $item[3]
    $ws.Cells.Item($row, 5).Value = "'" + $item[4]
}

Write-Output "Done updating cryptos data"